$d = $word.ActiveDocument

# 1. Merge the 3 runs in the "I think that everything is..." paragraph into one run
#    (the grammar-flagged "pretty straight" no longer needs special proofing runs).
$null = $d.Content.Find.Execute(
    "I think that everything is pretty straight forward this week and our jobs are clear.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "I think that everything is pretty straight forward this week and our jobs are clear.",
    2)

# 2. Merge the 3 runs in the "You have described and explored all of the
#    interactions..." paragraph into one run (same proofing cleanup).
$null = $d.Content.Find.Execute(
    "You have described and explored all of the interactions with other people.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "You have described and explored all of the interactions with other people.",
    2)

# 3. Move the _GoBack bookmark from right after "Jake" to the very start of the
#    document. Remove it from its old spot first.
$old = $d.Bookmarks.Item("_GoBack")
$old.Delete()

# Adding a bookmark directly at document position 0 mis-places the end tag in
# this runtime, so insert a throwaway character at the start, anchor the new
# bookmark right after it, then delete the throwaway character - the bookmark
# collapses back to a true (0,0) range in the first paragraph.
$startRange = $d.Range(0, 0)
$startRange.InsertBefore("X")
$bmRange = $d.Range(1, 1)
$d.Bookmarks.Add("_GoBack", $bmRange)
$d.Range(0, 1).Delete()
